$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data refresh timestamp (header cell A1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 01:03"

# Country rows re-ranked by updated "Casos totales" (column B); set every row
# A1:H218 data cell explicitly so both the country label in column A and the
# per-country stats in columns B:H land exactly where the refreshed, re-sorted
# dataset puts them.
$data = @(
    @("Estados Unidos", 1236141, 23306, 199151, 964823, 16174, 2246, 72167),
    @("España", 250561, 2260, 154718, 70230, 2254, 185, 25613),
    @("Italia", 213013, 1075, 85231, 98467, 1427, 236, 29315),
    @("Reino Unido", 194990, 4406, 0, 165219, 1559, 693, 29427),
    @("Francia", 170551, 1089, 52736, 92284, 3430, 330, 25531),
    @("Alemania", 167007, 855, 135100, 24914, 1937, 0, 6993),
    @("Rusia", 155370, 10102, 19865, 134054, 2300, 95, 1451),
    @("Turquia", 129491, 1832, 73285, 52686, 1338, 59, 3520),
    @("Brasil", 114715, 6449, 48221, 58573, 8318, 578, 7921),
    @("Iran", 99970, 1323, 80475, 13155, 2685, 63, 6340),
    @("China", 82881, 1, 77853, 395, 29, 0, 4633),
    @("Canada", 62038, 1266, 26938, 31061, 557, 185, 4039),
    @("Peru", 51189, 3817, 15413, 34332, 709, 100, 1444),
    @("Belgica", 50509, 242, 12441, 30052, 646, 92, 8016),
    @("India", 49400, 2963, 14142, 33565, 0, 127, 1693),
    @("Paises Bajos", 41087, 317, 0, 35669, 644, 86, 5168),
    @("Ecuador", 31881, 0, 3433, 26879, 159, 0, 1569),
    @("Arabia Saudita", 30251, 1595, 5431, 24620, 143, 9, 200),
    @("Suiza", 30009, 28, 25400, 2814, 121, 11, 1795),
    @("Portugal", 25702, 178, 1743, 22885, 134, 11, 1074),
    @("Mexico", 24905, 1434, 15938, 6696, 378, 117, 2271),
    @("Suecia", 23216, 495, 4074, 16288, 435, 85, 2854),
    @("Pakistan", 22049, 1108, 5801, 15734, 111, 38, 514),
    @("Chile", 22016, 1373, 10710, 11031, 470, 5, 275),
    @("Irlanda", 21983, 211, 13386, 7258, 93, 20, 1339),
    @("Singapur", 19410, 632, 1519, 17873, 24, 0, 18),
    @("Bielorrusia", 18350, 861, 3771, 14472, 92, 4, 107),
    @("Catar", 17142, 951, 1924, 15206, 72, 0, 12),
    @("Israel", 16289, 43, 10465, 5586, 89, 3, 238),
    @("Austria", 15650, 29, 13462, 1582, 104, 6, 606),
    @("Emiratos Arabes Unidos", 15192, 462, 3153, 11893, 1, 9, 146),
    @("Japon", 15078, 0, 4156, 10386, 328, 0, 536),
    @("Polonia", 14431, 425, 4280, 9435, 160, 18, 716),
    @("Rumania", 13837, 325, 5454, 7542, 244, 23, 841),
    @("Ucrania", 12697, 366, 1875, 10506, 167, 13, 316),
    @("Indonesia", 12071, 484, 2197, 9002, 0, 8, 872),
    @("Banglades", 10929, 786, 1403, 9343, 1, 1, 183),
    @("Corea del Sur", 10804, 3, 9283, 1267, 55, 2, 254),
    @("Dinamarca", 9821, 151, 7296, 2022, 49, 10, 503),
    @("Filipinas", 9684, 199, 1408, 7639, 31, 14, 637),
    @("Serbia", 9677, 120, 1723, 7754, 51, 3, 200),
    @("Colombia", 8613, 640, 2013, 6222, 120, 20, 378),
    @("Republica Dominicana", 8480, 245, 1905, 6221, 144, 8, 354),
    @("Noruega", 7955, 51, 32, 7708, 27, 1, 215),
    @("Chequia", 7896, 77, 4006, 3633, 59, 5, 257),
    @("Sudafrica", 7572, 352, 2746, 4678, 36, 10, 148),
    @("Panama", 7387, 190, 726, 6458, 93, 3, 203),
    @("Egipto", 7201, 388, 1730, 5019, 0, 16, 452),
    @("Australia", 6849, 24, 5889, 864, 27, 1, 96),
    @("Malasia", 6383, 30, 4567, 1710, 24, 1, 106),
    @("Kuwait", 5804, 526, 2032, 3732, 90, 0, 40),
    @("Finlandia", 5412, 85, 3500, 1666, 48, 6, 246),
    @("Marruecos", 5219, 166, 1838, 3200, 1, 2, 181),
    @("Argentina", 4887, 0, 1472, 3153, 148, 2, 262),
    @("Argelia", 4838, 190, 2067, 2301, 22, 5, 470),
    @("Moldavia", 4363, 115, 1544, 2683, 237, 4, 136),
    @("Kazajistan", 4205, 156, 1279, 2897, 40, 0, 29),
    @("Luxemburgo", 3840, 12, 3412, 332, 22, 0, 96),
    @("Barein", 3720, 187, 1762, 1950, 1, 0, 8),
    @("Afganistan", 3224, 330, 421, 2708, 7, 5, 95),
    @("Hungria", 3065, 30, 709, 1993, 55, 12, 363),
    @("Tailandia", 2988, 1, 2747, 187, 61, 0, 54),
    @("Nigeria", 2950, 148, 481, 2371, 4, 5, 98),
    @("Oman", 2735, 98, 858, 1865, 17, 0, 12),
    @("Ghana", 2719, 0, 294, 2407, 4, 0, 18),
    @("Grecia", 2642, 10, 1374, 1122, 35, 0, 146),
    @("Armenia", 2619, 112, 1111, 1468, 10, 1, 40),
    @("Irak", 2431, 85, 1571, 758, 0, 4, 102),
    @("Uzbekistan", 2207, 18, 1501, 696, 8, 0, 10),
    @("Croacia", 2112, 11, 1560, 469, 14, 3, 83),
    @("Camerun", 2104, 0, 953, 1087, 12, 0, 64),
    @("Azerbaiyan", 2060, 76, 1508, 526, 18, 0, 26),
    @("Bosnia y Herzegovina", 1946, 20, 911, 956, 4, 1, 79),
    @("Guinea", 1811, 101, 498, 1303, 0, 1, 10),
    @("Islandia", 1799, 0, 1733, 56, 0, 0, 10),
    @("Estonia", 1711, 8, 261, 1395, 6, 0, 55),
    @("Bulgaria", 1704, 52, 342, 1282, 37, 2, 80),
    @("Cuba", 1685, 17, 954, 662, 8, 0, 69),
    @("Bolivia", 1681, 87, 174, 1425, 3, 6, 82),
    @("Republica de Macedonia", 1526, 8, 1013, 427, 21, 1, 86),
    @("Nueva Zelanda", 1486, 0, 1302, 164, 0, 0, 20),
    @("Costa de Marfil", 1464, 32, 701, 745, 0, 1, 18),
    @("Eslovenia", 1445, 6, 244, 1103, 17, 1, 98),
    @("Lituania", 1423, 4, 678, 699, 17, 0, 46),
    @("Eslovaquia", 1421, 8, 741, 655, 6, 0, 25),
    @("Senegal", 1329, 58, 470, 848, 6, 1, 11),
    @("Honduras", 1178, 123, 122, 973, 10, 1, 83),
    @("Republica de Yibuti", 1120, 4, 745, 373, 0, 0, 2),
    @("Hong Kong", 1041, 0, 920, 117, 1, 0, 4),
    @("Tunez", 1022, 4, 482, 497, 17, 0, 43),
    @("Letonia", 896, 0, 348, 531, 4, 1, 17),
    @("Republica de Chipre", 878, 4, 296, 567, 15, 0, 15),
    @("Kirguistan", 843, 13, 600, 232, 13, 1, 11),
    @("Somalia", 835, 79, 75, 722, 2, 3, 38),
    @("Albania", 820, 17, 570, 219, 7, 0, 31),
    @("Sudan", 778, 100, 70, 663, 0, 4, 45),
    @("Sri Lanka", 771, 20, 213, 549, 1, 1, 9),
    @("Niger", 763, 8, 543, 182, 0, 1, 38),
    @("Costa Rica", 755, 13, 413, 336, 5, 0, 6),
    @("Principado de Andorra", 751, 1, 514, 191, 16, 1, 46),
    @("Libano", 741, 1, 206, 510, 43, 0, 25),
    @("Mayotte", 739, 53, 352, 378, 6, 3, 9),
    @("Guatemala", 730, 27, 79, 632, 5, 2, 19),
    @("Crucero", 712, 0, 645, 54, 4, 0, 13),
    @("Consejo Danes para los Refugiados", 705, 23, 90, 581, 0, 0, 34),
    @("Burkina Faso", 688, 16, 548, 92, 0, 2, 48),
    @("Uruguay", 657, 0, 447, 193, 10, 0, 17),
    @("Mali", 612, 32, 228, 352, 0, 3, 32),
    @("Georgia", 604, 11, 240, 355, 6, 0, 9),
    @("San Marino", 589, 7, 92, 456, 5, 0, 41),
    @("El Salvador", 587, 32, 201, 373, 3, 0, 13),
    @("Maldivas", 573, 32, 20, 551, 2, 1, 2),
    @("Kenia", 535, 45, 182, 329, 2, 0, 24),
    @("Malta", 482, 2, 403, 74, 1, 1, 5),
    @("Tanzania", 480, 0, 167, 297, 7, 0, 16),
    @("Jamaica", 471, 2, 49, 413, 2, 0, 9),
    @("Jordania", 471, 6, 375, 87, 5, 0, 9),
    @("Taiwan", 438, 0, 334, 98, 0, 0, 6),
    @("Paraguay", 431, 35, 135, 286, 7, 0, 10),
    @("Reunion", 424, 0, 300, 124, 2, 0, 0),
    @("Guinea-Bisau", 413, 0, 19, 393, 0, 0, 1),
    @("Gabon", 397, 30, 93, 298, 1, 0, 6),
    @("Estado de Palestina", 371, 9, 127, 242, 0, 0, 2),
    @("Venezuela", 361, 4, 158, 193, 1, 0, 10),
    @("Mauricio", 332, 0, 319, 3, 3, 0, 10),
    @("Isla de Man", 326, 1, 271, 32, 19, 0, 23),
    @("Montenegro", 324, 1, 253, 63, 2, 0, 8),
    @("Guinea Ecuatorial", 315, 0, 13, 299, 0, 0, 3),
    @("Tayikistan", 293, 63, 0, 288, 0, 2, 5),
    @("Vietnam", 271, 0, 232, 39, 8, 0, 0),
    @("Ruanda", 261, 0, 129, 132, 0, 0, 0),
    @("Congo", 236, 0, 26, 200, 0, 0, 10),
    @("Sierra Leona", 199, 21, 43, 145, 0, 2, 11),
    @("Islas Feroe", 187, 0, 185, 2, 0, 0, 0),
    @("Cabo Verde", 186, 11, 37, 147, 0, 0, 2),
    @("Martinica", 181, 0, 83, 84, 5, 0, 14),
    @("Santo Tome y Principe", 174, 151, 4, 167, 0, 0, 3),
    @("Republica del Chad", 170, 53, 43, 110, 0, 7, 17),
    @("Liberia", 170, 4, 58, 92, 0, 2, 20),
    @("Birmania", 161, 0, 49, 106, 0, 0, 6),
    @("Guadalupe", 152, 0, 98, 42, 5, 0, 12),
    @("Madagascar", 151, 2, 101, 50, 1, 0, 0),
    @("Etiopia", 145, 5, 91, 50, 0, 1, 4),
    @("Gibraltar", 144, 0, 133, 11, 0, 0, 0),
    @("Zambia", 138, 1, 92, 43, 1, 0, 3),
    @("Brunei", 138, 0, 131, 6, 2, 0, 1),
    @("Guayana Francesa", 133, 0, 100, 32, 2, 0, 1),
    @("Togo", 128, 2, 74, 45, 0, 0, 9),
    @("Camboya", 122, 0, 120, 2, 1, 0, 0),
    @("Suazilandia", 119, 3, 12, 106, 0, 0, 1),
    @("Trinidad yTobago", 116, 0, 102, 6, 0, 0, 8),
    @("Bermudas", 115, 0, 54, 54, 4, 0, 7),
    @("Aruba", 101, 1, 82, 17, 4, 0, 2),
    @("Haiti", 100, 0, 10, 79, 0, 0, 11),
    @("Uganda", 98, 1, 55, 43, 0, 0, 0),
    @("Benin", 96, 0, 50, 44, 0, 0, 2),
    @("Monaco", 95, 0, 81, 10, 1, 0, 4),
    @("Guyana", 92, 0, 27, 56, 3, 0, 9),
    @("Bahamas", 89, 6, 26, 52, 1, 0, 11),
    @("Republica de Africa Central", 85, 0, 10, 75, 0, 0, 0),
    @("Nepal", 82, 7, 16, 66, 0, 0, 0),
    @("Barbados", 82, 0, 47, 28, 4, 0, 7),
    @("Liechtenstein", 82, 0, 55, 26, 0, 0, 1),
    @("Mozambique", 81, 1, 19, 62, 0, 0, 0),
    @("San Martin (Parte Holandesa)", 76, 0, 44, 18, 7, 1, 14),
    @("Islas Caimanes", 75, 0, 14, 60, 3, 0, 1),
    @("Libia", 63, 0, 23, 37, 0, 0, 3),
    @("Polinesia Francesa", 58, 0, 53, 5, 1, 0, 0),
    @("Sudan del Sur", 52, 6, 0, 52, 0, 0, 0),
    @("Macao", 45, 0, 39, 6, 1, 0, 0),
    @("Siria", 44, 0, 27, 14, 0, 0, 3),
    @("Malaui", 41, 0, 9, 29, 1, 0, 3),
    @("Mongolia", 41, 1, 13, 28, 0, 0, 0),
    @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
    @("Eritrea", 39, 0, 30, 9, 0, 0, 0),
    @("San Martin (Parte Francesa)", 38, 0, 29, 6, 1, 0, 3),
    @("Angola", 36, 1, 11, 23, 0, 0, 2),
    @("Zimbabue", 34, 0, 5, 25, 0, 0, 4),
    @("Guam", 32, 0, 0, 31, 0, 0, 1),
    @("Antigua y Barbuda", 25, 0, 16, 6, 1, 0, 3),
    @("Timor Oriental", 24, 0, 20, 4, 0, 0, 0),
    @("Botsuana", 23, 0, 8, 14, 0, 0, 1),
    @("Yemen", 22, 10, 1, 17, 0, 2, 4),
    @("Granada", 21, 0, 13, 8, 4, 0, 0),
    @("Laos", 19, 0, 9, 10, 0, 0, 0),
    @("Fiyi", 18, 0, 14, 4, 0, 0, 0),
    @("Santa Lucia", 18, 0, 15, 3, 0, 0, 0),
    @("Nueva Caledonia", 18, 0, 17, 1, 1, 0, 0),
    @("Belice", 18, 0, 16, 0, 0, 0, 2),
    @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @("San Vicente y las Granadinas", 17, 0, 9, 8, 0, 0, 0),
    @("Gambia", 17, 0, 9, 7, 0, 0, 1),
    @("Namibia", 16, 0, 8, 8, 0, 0, 0),
    @("Nicaragua", 16, 1, 7, 4, 0, 0, 5),
    @("Dominica", 16, 0, 14, 2, 0, 0, 0),
    @("Curazao", 16, 0, 13, 2, 0, 0, 1),
    @("Burundi", 15, 0, 7, 7, 0, 0, 1),
    @("San Cristobal y Nieves", 15, 0, 12, 3, 0, 0, 0),
    @("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0),
    @("Islas Turcas y Caicos", 12, 0, 6, 5, 0, 0, 1),
    @("Santa Sede", 11, 0, 2, 9, 0, 0, 0),
    @("Montserrat", 11, 0, 7, 3, 1, 0, 1),
    @("Seychelles", 11, 0, 8, 3, 0, 0, 0),
    @("Groenlandia", 11, 0, 11, 0, 0, 0, 0),
    @("Surinam", 10, 0, 9, 0, 0, 0, 1),
    @("Mauritania", 8, 0, 6, 1, 0, 0, 1),
    @("Papua Nueva Guinea", 8, 0, 8, 0, 0, 0, 0),
    @("Butan", 7, 0, 5, 2, 0, 0, 0),
    @("Bonaire, San Eustaquio y Saba", 6, 0, 0, 6, 0, 0, 0),
    @("Islas Virgenes Britanicas", 6, 0, 3, 2, 0, 0, 1),
    @("Sahara Occidental", 6, 0, 5, 1, 0, 0, 0),
    @("San Bartolome", 6, 0, 6, 0, 0, 0, 0),
    @("Comoras", 3, 0, 0, 3, 0, 0, 0),
    @("Anguila", 3, 0, 3, 0, 0, 0, 0),
    @("San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

